# Añadidos los primeros comandos, la clase para encriptar, eliminadas
# algunas cosas y añadidas otras.
#
# Adds a new citizen record (row 5) to Hoja1, mirroring the layout of the
# existing rows 2-4: Nombre, Apellidos, Correo electronico (as a mailto
# hyperlink), Fecha nacimiento (date-formatted), Direccion postal,
# Nacionalidad and DNI. Also moves the active selection, matching the
# author's last interaction with the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row of data -------------------------------------------------
$ws.Range("A5").Value = "Seila"
$ws.Range("B5").Value = "Khayat Prada"
$ws.Range("C5").Value = "porqueestonova@joder.com"
$ws.Range("D5").Value = 35079
$ws.Range("E5").Value = "Deberia Formatear :)"
$ws.Range("F5").Value = "Español"
$ws.Range("G5").Value = "34234239P"

# --- Hyperlink the e-mail cell, like the rows above it ---------------
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:porqueestonova@joder.com")

# --- Match formatting of the existing rows ----------------------------
# Email column uses the hyperlink style.
$ws.Range("C2").Copy()
$ws.Range("C5").PasteSpecial(-4122)

# Birth-date column uses a date number format.
$ws.Range("D2").Copy()
$ws.Range("D5").PasteSpecial(-4122)

# DNI column picks up the formatting used on row 4.
$ws.Range("G4").Copy()
$ws.Range("G5").PasteSpecial(-4122)

# --- Selection left where the author last clicked ----------------------
$ws.Range("D10").Select()
